# Insert a new weekly price record as row 358 ("Fruta / hortaliza, semanal").
# This shifts the former rows 358-432 down to 359-433 (unchanged otherwise)
# and fills the freshly inserted row 358 with the new observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at 358, pushing existing rows 358:432 down to 359:433.
$ws.Rows("358:358").Insert()

$ws.Range("A358").Value = 6
$ws.Range("B358").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C358").Value = "Metropolitana"
$ws.Range("D358").Value = 44694
$ws.Range("E358").Value = 13
$ws.Range("F358").Value = 100112039
$ws.Range("G358").Value = "Ciboulette"
$ws.Range("H358").Value = "Sin especificar"
$ws.Range("I358").Value = "Primera"
$ws.Range("J358").Value = 740
$ws.Range("K358").Value = 700
$ws.Range("L358").Value = 800
$ws.Range("M358").Value = 747
$ws.Range("N358").Value = "`$/docena de atados"
$ws.Range("O358").Value = "Región Metropolitana"
$ws.Range("P358").Value = 249
$ws.Range("Q358").Value = 3
$ws.Range("R358").Value = "Hortaliza"
